# EnergiaCore.xlsx update
# - Files sheet: remove twi.c/twi.h/usci_isr_handler.c/usci_isr_handler.h/Wire.cpp/Wire.h rows
#                and update Status (col A) markers for several files
# - Functions sheet: mark several previously-unmarked functions as implemented (col C = "X")
# - Update the remembered selection on both sheets

$wb  = $excel.ActiveWorkbook
$wsFiles     = $wb.Worksheets.Item("Files")
$wsFunctions = $wb.Worksheets.Item("Functions")

# ---------------------------------------------------------------------------
# 1. "Files" sheet: drop the rows for files that are no longer tracked here
#    (twi.c, twi.h, usci_isr_handler.c, usci_isr_handler.h, Wire.cpp, Wire.h).
#    Delete from the bottom up so earlier row numbers stay valid.
# ---------------------------------------------------------------------------
$wsFiles.Rows(26).Delete()   # Wire.h
$wsFiles.Rows(25).Delete()   # Wire.cpp
$wsFiles.Rows(22).Delete()   # usci_isr_handler.h
$wsFiles.Rows(21).Delete()   # usci_isr_handler.c
$wsFiles.Rows(20).Delete()   # twi.h
$wsFiles.Rows(19).Delete()   # twi.c

# ---------------------------------------------------------------------------
# 2. "Files" sheet: update the Status column for the remaining rows
# ---------------------------------------------------------------------------
$wsFiles.Range("A4").Value2  = $null   # Energia.h        "/" -> (blank)
$wsFiles.Range("A5").Value2  = "X"     # HardwareSerial.cpp (blank) -> X
$wsFiles.Range("A6").Value2  = "X"     # HardwareSerial.h   (blank) -> X
$wsFiles.Range("A7").Value2  = "X"     # lm4f.h            "/" -> X
$wsFiles.Range("A18").Value2 = $null   # Tone.cpp           X  -> (blank)
$wsFiles.Range("A24").Value2 = "X"     # wiring_private.c  (blank) -> X
$wsFiles.Range("A26").Value2 = $null   # wiring_shift.c     X  -> (blank)
$wsFiles.Range("A30").Value2 = $null   # pins_energia.h    "/" -> (blank)

# ---------------------------------------------------------------------------
# 3. "Functions" sheet: mark additional functions as done (column C)
# ---------------------------------------------------------------------------
$wsFunctions.Range("C20").Value2 = "X"   # randomSeed()
$wsFunctions.Range("C21").Value2 = "X"   # random()

38..53 | ForEach-Object {
    $wsFunctions.Range("C$_").Value2 = "X"
}

# ---------------------------------------------------------------------------
# 4. Restore the remembered cell selections (select Files' cell first so the
#    Functions sheet - the real active tab - ends up selected last).
# ---------------------------------------------------------------------------
$wsFiles.Range("A8").Select()
$wsFunctions.Range("C22").Select()
